$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update rows 20 and 21: Priority -> 10, Test Done -> "yes", Name -> "transcript",
# comments -> "just regtest"
# (comments are updated first, on both rows, so the shared string that used to
# hold "retest, look up in refsew, ensembl" is freed up and reused in place for
# "just regtest" before the brand-new "transcript" string is appended.)
$ws.Range("A20").Value = 10
$ws.Range("C20").Value = "yes"
$ws.Range("E20").Value = "just regtest"

$ws.Range("A21").Value = 10
$ws.Range("C21").Value = "yes"
$ws.Range("E21").Value = "just regtest"

$ws.Range("D20").Value = "transcript"
$ws.Range("D21").Value = "transcript"

# Update the view: scroll so that A7 is the top-left visible cell, and select B21
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("B21").Select()
